{"js": "// 1) \"Event organizers ...\" -> \"Co-creation event organizers ...\"\n//    (the capital \"E\" of \"Event\" becomes part of \"Co-creation e\", then \"vent organizers...\")\nconst introResults = context.document.body.search(\"Event\", { matchCase: true });\nintroResults.load(\"text\");\nawait context.sync();\nintroResults.items[0].insertText(\"Co-creation event\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 2) The \"_GoBack\" bookmark moves from the empty paragraph near the end of the\n//    document to the end of the paragraph we just edited (\"... visit multiple\n//    websites.\"). Remove it from its old location first.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 3) Re-locate the edited paragraph (index may have shifted) and drop the\n//    bookmark at its very end.\nconst bodyParagraphs = context.document.body.paragraphs;\nbodyParagraphs.load(\"text\");\nawait context.sync();\n\nlet introParaIndex = -1;\nfor (let i = 0; i < bodyParagraphs.items.length; i++) {\n  if (bodyParagraphs.items[i].text.indexOf(\"visit multiple websites\") >= 0) {\n    introParaIndex = i;\n    break;\n  }\n}\nconst introPara = bodyParagraphs.items[introParaIndex];\nconst introParaEnd = introPara.getRange(\"End\");\nintroParaEnd.insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// 4) \"The concept of the semantic web is making websites reada\" + \"ble by\n//    machines. \" were stored as two runs with identical formatting; they are\n//    now a single run. Re-writing the paragraph's own text collapses the\n//    (formatting-identical) runs into one without altering the content.\nbodyParagraphs.load(\"text\");\nawait context.sync();\n\nlet readableParaIndex = -1;\nfor (let i = 0; i < bodyParagraphs.items.length; i++) {\n  if (bodyParagraphs.items[i].text.indexOf(\"making websites reada\") >= 0) {\n    readableParaIndex = i;\n    break;\n  }\n}\nconst readablePara = bodyParagraphs.items[readableParaIndex];\nreadablePara.load(\"text\");\nawait context.sync();\nreadablePara.getRange().insertText(readablePara.text, Word.InsertLocation.replace);\nawait context.sync();\n\n// 5) Likewise for \"Events are always associated ...\" + \" The semantic web\n//    allows ...\" which merge into a single run.\nbodyParagraphs.load(\"text\");\nawait context.sync();\n\nlet eventsDataParaIndex = -1;\nfor (let i = 0; i < bodyParagraphs.items.length; i++) {\n  if (bodyParagraphs.items[i].text.indexOf(\"Events are always associated\") >= 0) {\n    eventsDataParaIndex = i;\n    break;\n  }\n}\nconst eventsDataPara = bodyParagraphs.items[eventsDataParaIndex];\neventsDataPara.load(\"text\");\nawait context.sync();\neventsDataPara.getRange().insertText(eventsDataPara.text, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) \"Event organizers ...\" -> \"Co-creation event organizers ...\"\n#    (the capital \"E\" of \"Event\" becomes part of \"Co-creation e\", then \"vent organizers...\")\n$introRng = $d.Content\n$introRng.Find.Execute(\"Event\", $false, $true)\n$introRng.Text = \"Co-creation event\"\n\n# 2) The \"_GoBack\" bookmark moves from the empty paragraph near the end of the\n#    document to the end of the paragraph we just edited (\"... visit multiple\n#    websites.\"). Remove it from its old location first.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# 3) Re-locate the edited paragraph (paragraph 4, \"Co-creation event organizers\n#    ... visit multiple websites.\") and drop the bookmark at its very end,\n#    still inside the paragraph (i.e. before the paragraph mark).\n$introPara = $d.Paragraphs.Item(4)\n$introEnd = $introPara.Range\n$introEnd.MoveEnd(1, -1)\n$introEnd.Collapse(0)\n$introEnd.InsertAfter(\"TmpGoBackMarker\")\n$d.Bookmarks.Add(\"_GoBack\", $introEnd)\n$markerRng = $d.Content\n$markerRng.Find.Execute(\"TmpGoBackMarker\")\n$markerRng.Text = \"\"\n\n# 4) \"The concept of the semantic web is making websites reada\" + \"ble by\n#    machines. \" were stored as two runs with identical formatting; they are\n#    now a single run. Temporarily changing and restoring the paragraph's own\n#    text forces the engine to rebuild the run, collapsing the\n#    (formatting-identical) runs into one without altering the content.\n$readablePara = $d.Paragraphs.Item(11)\n$readableRng = $readablePara.Range\n$readableRng.MoveEnd(1, -1)\n$readableOriginal = $readableRng.Text\n$readableRng.Text = $readableOriginal + \"ZtmpZ\"\n$readablePara2 = $d.Paragraphs.Item(11)\n$readableRng2 = $readablePara2.Range\n$readableRng2.MoveEnd(1, -1)\n$readableRng2.Text = $readableOriginal\n\n# 5) Likewise for \"Events are always associated ...\" + \" The semantic web\n#    allows ...\" which merge into a single run.\n$eventsDataPara = $d.Paragraphs.Item(13)\n$eventsDataRng = $eventsDataPara.Range\n$eventsDataRng.MoveEnd(1, -1)\n$eventsDataOriginal = $eventsDataRng.Text\n$eventsDataRng.Text = $eventsDataOriginal + \"ZtmpZ\"\n$eventsDataPara2 = $d.Paragraphs.Item(13)\n$eventsDataRng2 = $eventsDataPara2.Range\n$eventsDataRng2.MoveEnd(1, -1)\n$eventsDataRng2.Text = $eventsDataOriginal\n"}
